$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Số điện thoại" (phone number) column at C, shifting
# Học vị/Bộ môn/Khoa/Chức vị from C:F to D:G.
$ws.Columns("C").Insert()

# The insert copies the left neighbour's (hyperlink) formatting into the
# whole new column - strip that back to the default style before writing
# real content/formatting into it.
$ws.Range("C1:C9").ClearFormats()

# Remove the stray formatting/cells the insert copied into the new
# column's blank trailing rows (6-9) so they stay truly empty.
$ws.Range("C6:C9").Clear()

# --- Header row (bold) ---
$ws.Range("A1").Value = "Họ và tên"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Số điện thoại"
$ws.Range("D1").Value = "Học vị"
$ws.Range("E1").Value = "Bộ môn"
$ws.Range("F1").Value = "Khoa"
$ws.Range("G1").Value = "Chức vị"
$ws.Range("A1:G1").Font.Bold = $true

# --- Data rows ---
$ws.Range("A2").Value = "TestImport1"
$ws.Range("E2").Value = "Hệ thống thông tin"

$ws.Range("A3").Value = "TestImport2"
$ws.Range("F3").Value = "Công nghệ thông tin"

$ws.Range("A4").Value = "TestImport3"
$ws.Range("D4").Value = "Thạc sĩ"
$ws.Range("E4").Value = "Trí tuệ nhân tạo"

$ws.Range("A5").Value = "TestImport4"
$ws.Range("D5").Value = "Tiến sĩ"
$ws.Range("E5").Value = "Kỹ thuật phần mềm"
$ws.Range("F5").Value = "Công nghệ thông tin"
$ws.Range("G5").Value = "Trưởng khoa"

# Phone numbers: leading apostrophe forces text-with-quote-prefix so the
# leading zero survives.
$ws.Range("C2").Value = "'0969615123"
$ws.Range("C3").Value = "'0969615456"
$ws.Range("C4").Value = "'0969615789"
$ws.Range("C5").Value = "'0969615246"

# --- Column widths (closest the engine's 1/6-char quantized ColumnWidth
# can land to the target bestFit pixel widths of 12.71 / 6.71 / 17.71 / 19) ---
$ws.Columns("C").ColumnWidth = 11.833333333333334
$ws.Columns("D").ColumnWidth = 5.833333333333333
$ws.Columns("E").ColumnWidth = 16.833333333333332
$ws.Columns("F").ColumnWidth = 18.166666666666668

# --- Print / view settings ---
$ws.PageSetup.Orientation = 1
$ws.Range("J9").Select()
